$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2 = 34.1502
    3 = 26.359
    4 = 23.2081
    5 = 21.6382
    6 = 20.2516
    7 = 19.251
    8 = 18.6105
    9 = 18.0149
    10 = 17.4068
    11 = 16.9958
    12 = 16.4747
    13 = 16.1425
    14 = 15.9092
    15 = 15.3756
    16 = 15.1456
    17 = 14.5205
    18 = 14.1912
    19 = 13.7762
    20 = 13.3218
    21 = 12.9442
    22 = 12.4055
    23 = 11.902
    24 = 11.4914
    25 = 11.1672
    26 = 10.8755
    27 = 10.5402
    28 = 10.4516
    29 = 10.2334
    30 = 10.0067
    31 = 9.890700000000001
    32 = 9.738
    33 = 9.543699999999999
    34 = 9.3931
    35 = 9.3147
    36 = 9.333600000000001
    37 = 9.205500000000001
    38 = 9.1509
    39 = 9.017200000000001
    40 = 9.0693
    41 = 8.8813
    42 = 8.9475
    43 = 8.788
    44 = 8.818300000000001
    45 = 8.7014
    46 = 8.722
    47 = 8.642300000000001
    48 = 8.5817
    49 = 8.488099999999999
    50 = 8.644399999999999
    51 = 8.428900000000001
    52 = 8.393599999999999
    53 = 8.4091
    54 = 8.416700000000001
    55 = 8.2903
    56 = 8.2728
    57 = 8.278
    58 = 8.146000000000001
    59 = 8.1754
    60 = 8.048500000000001
    61 = 8.1318
    62 = 8.139699999999999
    63 = 8.1144
    64 = 7.9396
    65 = 7.9627
    66 = 7.9509
    67 = 7.8626
    68 = 7.8828
    69 = 7.9132
    70 = 8.0014
    71 = 7.9219
    72 = 7.9201
    73 = 7.8476
    74 = 7.838
    75 = 7.8439
    76 = 7.7219
    77 = 7.7787
    78 = 7.8238
    79 = 7.8293
    80 = 7.7649
    81 = 7.7102
    82 = 7.6923
    83 = 7.7038
    84 = 7.6921
    85 = 7.6737
}

foreach ($r in $values.Keys) {
    $ws.Cells.Item($r, 4).Value = $values[$r]
}